$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Defects Found")

$ws.Range("A5").Value = "Route Calculations"
$ws.Range("B5").Value = "Route calculated is not the shortest possible path"
$ws.Range("C5").Value = 3.3
$ws.Range("D5").Value = "Correctness"

$ws.Range("A5").WrapText = $true
